# Update existing row 2 values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: "asa" -> "santiago"
$ws.Range("A2").Value = "santiago"

# G2: "07-01-2015" -> "02-01-2025" (keep as literal text, not a date)
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "02-01-2025"
$ws.Range("G2").ClearFormats()

# Append new row 4 with the new applicant's data
$ws.Range("A4").Value = "SANTIAGO"
$ws.Range("B4").Value = "RAMIREZ"
$ws.Range("C4").Value = "VALENCIA"
$ws.Range("D4").Value = "C.C."

# E4 looks numeric, force text so it stays "1001456789"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1001456789"
$ws.Range("E4").ClearFormats()

$ws.Range("F4").Value = "CRA 64 C NRO 103 - 41"

# G4 looks like a date, force text so it stays "01-09-2002"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "01-09-2002"
$ws.Range("G4").ClearFormats()

$ws.Range("H4").Value = 22
$ws.Range("I4").Value = "GIRARDOT"
$ws.Range("J4").Value = 3

# K4 looks numeric, force text so it stays "1000121514"
$ws.Range("K4").NumberFormat = "@"
$ws.Range("K4").Value = "1000121514"
$ws.Range("K4").ClearFormats()

$ws.Range("L4").Value = "santiago.@gmail.com"
$ws.Range("M4").Value = "ADMINISTRACIÓN DE EMPRESAS"
$ws.Range("N4").Value = 5
$ws.Range("O4").Value = "uploaded_files\1001456789_CÉDULA.pdf"
$ws.Range("P4").Value = "uploaded_files\1001456789_CIVICA.pdf"
$ws.Range("Q4").Value = "uploaded_files\1001456789_SERVICIOPUBLICOS.pdf"
$ws.Range("R4").Value = "uploaded_files\1001456789_ANEXO1.pdf"
$ws.Range("S4").Value = "uploaded_files\1001456789_ANEXO2.pdf"
